# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for Naranja (Macroferia Regional de Talca)
# right before the current row 269, shifting the existing data (old rows
# 269..338) down to rows 272..341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 269:271, shifting everything below down (xlShiftDown = -4121)
$ws.Rows("269:271").Insert(-4121)

# ---- New row 269: Naranja / Lane Late / Primera ----
$ws.Range("A269").Value = 5
$ws.Range("B269").Value = "Macroferia Regional de Talca"
$ws.Range("C269").Value = "Maule"
$ws.Range("D269").Value = 44508
$ws.Range("E269").Value = 7
$ws.Range("F269").Value = "Fruta"
$ws.Range("G269").Value = 100102
$ws.Range("H269").Value = "Cítricos"
$ws.Range("I269").Value = 100102005
$ws.Range("J269").Value = "Naranja"
$ws.Range("K269").Value = "Lane Late"
$ws.Range("L269").Value = "Primera"
$ws.Range("M269").Value = 250
$ws.Range("N269").Value = 8000
$ws.Range("O269").Value = 8000
$ws.Range("P269").Value = 8000
$ws.Range("Q269").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R269").Value = "Región de O'Higgins"
$ws.Range("S269").Value = 533
$ws.Range("T269").Value = 15

# ---- New row 270: Naranja / Navel Late / Primera ----
$ws.Range("A270").Value = 5
$ws.Range("B270").Value = "Macroferia Regional de Talca"
$ws.Range("C270").Value = "Maule"
$ws.Range("D270").Value = 44508
$ws.Range("E270").Value = 7
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100102
$ws.Range("H270").Value = "Cítricos"
$ws.Range("I270").Value = 100102005
$ws.Range("J270").Value = "Naranja"
$ws.Range("K270").Value = "Navel Late"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 300
$ws.Range("N270").Value = 8000
$ws.Range("O270").Value = 8000
$ws.Range("P270").Value = 8000
$ws.Range("Q270").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R270").Value = "Provincia de Melipilla"
$ws.Range("S270").Value = 533
$ws.Range("T270").Value = 15

# ---- New row 271: Naranja / Olinda Valencia / Primera ----
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value = 44508
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = "Fruta"
$ws.Range("G271").Value = 100102
$ws.Range("H271").Value = "Cítricos"
$ws.Range("I271").Value = 100102005
$ws.Range("J271").Value = "Naranja"
$ws.Range("K271").Value = "Olinda Valencia"
$ws.Range("L271").Value = "Primera"
$ws.Range("M271").Value = 400
$ws.Range("N271").Value = 8000
$ws.Range("O271").Value = 8000
$ws.Range("P271").Value = 8000
$ws.Range("Q271").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R271").Value = "Región de O'Higgins"
$ws.Range("S271").Value = 533
$ws.Range("T271").Value = 15
